$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.326106071472168
$ws.Range("B1").Value = 5.890065670013428
$ws.Range("C1").Value = 7.081566333770752
$ws.Range("D1").Value = 9.575989723205566
$ws.Range("E1").Value = 4.582254886627197
